# Generate Report for Archive
#
# The localization status for the outstanding items moved on from the
# handoff stage, so the workbook's status text needs updating everywhere
# it is shown (the Overview roll-up sheet plus each per-language detail
# sheet), and the now-narrower status columns are re-sized to fit the
# shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# "Ready for handoff" -> "In Translation" everywhere it is used.
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# The status columns can now be narrower to match the new text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
